{"js": "// The document contains a Word \"complex field\" (fldChar begin / instrText\n// \"m:self.name\" / fldChar separate / ... / fldChar end) that encodes an\n// M2Doc query. The edit rewrites that field into plain text runs holding\n// the M2Doc token syntax \"{m:self.name}\", split across four runs:\n//   \"{m:self.name\"  |  \"\"  |  \"\"  |  \"}\"\n// (mirrors TokenIteratorFieldRewriterSplit, which splits the rewritten\n// token text across the same number of runs the field used to occupy).\n\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\nif (fields.items.length === 0) {\n  throw new Error(\"No fields found in the document body.\");\n}\n\n// Locate the field carrying the \"m:self.name\" M2Doc query.\nfields.load(\"items/code\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < fields.items.length; i++) {\n  const code = (fields.items[i].code || \"\").trim();\n  if (code === \"m:self.name\") {\n    target = fields.items[i];\n    break;\n  }\n}\nif (!target) {\n  target = fields.items[0];\n}\n\n// Remember which paragraph holds the field so we can re-insert the\n// rewritten text runs at the same spot once the field is gone.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet hostParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  // The field paragraph has no visible text (fldChar/instrText runs do not\n  // surface through .text), so an empty paragraph is our field's host.\n  if (paragraphs.items[i].text === \"\") {\n    hostParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!hostParagraph) {\n  throw new Error(\"Could not locate the paragraph hosting the field.\");\n}\n\n// Remove the field itself (drops the fldChar/instrText runs it owns).\ntarget.delete();\nawait context.sync();\n\n// Re-insert the M2Doc token text as four plain-text runs, matching the\n// rewriter's output: \"{m:self.name\" + \"\" + \"\" + \"}\".\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:rPr/><w:t>{m:self.name</w:t></w:r>\n            <w:r><w:rPr/><w:t/></w:r>\n            <w:r><w:rPr/><w:t/></w:r>\n            <w:r><w:rPr/><w:t>}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nhostParagraph.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "# The document contains a Word \"complex field\" (fldChar begin / instrText\n# \"m:self.name\" / fldChar separate / ... / fldChar end) that encodes an\n# M2Doc query. This rewrites that field into plain text runs holding the\n# M2Doc token syntax \"{m:self.name}\", split across four runs:\n#   \"{m:self.name\"  |  \"\"  |  \"\"  |  \"}\"\n# (mirrors TokenIteratorFieldRewriterSplit, which splits the rewritten\n# token text across the same number of runs the field used to occupy.)\n\n$d = $word.ActiveDocument\n\n# Locate the field carrying the \"m:self.name\" M2Doc query.\n$targetField = $null\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n  $fld = $d.Fields.Item($i)\n  $code = $fld.Code.Text.Trim()\n  if ($code -eq \"m:self.name\") {\n    $targetField = $fld\n    break\n  }\n}\nif ($targetField -eq $null -and $d.Fields.Count -gt 0) {\n  $targetField = $d.Fields.Item(1)\n}\nif ($targetField -eq $null) {\n  throw \"No field found to rewrite.\"\n}\n\n# Find the paragraph that hosts the field so we replace its content\n# in-place (keeps that paragraph's own pPr/pStyle formatting).\n$fieldStart = $targetField.Code.Start\n$hostPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($fieldStart -ge $p.Range.Start -and $fieldStart -lt $p.Range.End) {\n    $hostPara = $p\n    break\n  }\n}\nif ($hostPara -eq $null) {\n  throw \"Could not locate the paragraph hosting the field.\"\n}\n\n$r = $hostPara.Range\n\n# Replace the paragraph's content (the field) with four plain-text runs\n# carrying the M2Doc token syntax: \"{m:self.name\" + \"\" + \"\" + \"}\".\n$ooxml = '<?xml version=\"1.0\" standalone=\"yes\"?><?mso-application progid=\"Word.Document\"?>' + `\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n'<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"512\"><pkg:xmlData>' + `\n'<Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\"><Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/></Relationships>' + `\n'</pkg:xmlData></pkg:part>' + `\n'<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData>' + `\n'<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"Normal\"/><w:rPr/></w:pPr><w:r><w:rPr/><w:t>{m:self.name</w:t></w:r><w:r><w:rPr/><w:t/></w:r><w:r><w:rPr/><w:t/></w:r><w:r><w:rPr/><w:t>}</w:t></w:r></w:p></w:body></w:document>' + `\n'</pkg:xmlData></pkg:part>' + `\n'</pkg:package>'\n\n$r.InsertXML($ooxml)\n"}
